$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 12500
$ws.Range("I43").Value = 12500
$ws.Range("K43").Value = 12500
$ws.Range("M43").Value = -12431

# Row 55
$ws.Range("H55").Value = 834.55554
$ws.Range("I55").Value = 469.30768
$ws.Range("J55").Value = 1784.2
$ws.Range("K55").Value = 469.30768
$ws.Range("L55").Value = 1784.2
$ws.Range("M55").Value = -255.30768
$ws.Range("N55").Value = -2212.2

# Row 100
$ws.Range("H100").Value = 2342.5386
$ws.Range("I100").Value = 1828.5555
$ws.Range("J100").Value = 3499
$ws.Range("K100").Value = 1828.5555
$ws.Range("L100").Value = 3499
$ws.Range("M100").Value = -1287.5555
$ws.Range("N100").Value = -4581

# Row 131
$ws.Range("H131").Value = 1290.7693
$ws.Range("I131").Value = 868
$ws.Range("K131").Value = 2604
$ws.Range("M131").Value = 2436

# Row 138
$ws.Range("H138").Value = 4781.655
$ws.Range("I138").Value = 3634.9312
$ws.Range("J138").Value = 6060.6924
$ws.Range("K138").Value = 10904.7936
$ws.Range("L138").Value = 18182.0772
$ws.Range("M138").Value = -5764.793600000001
$ws.Range("N138").Value = -28462.0772

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 1855.1428
$ws.Range("J88").Value = 1997.25
$ws.Range("L88").Value = 1997.25
$ws.Range("N88").Value = -2809.25

# Row 91
$ws.Range("H91").Value = 1855.1428
$ws.Range("J91").Value = 1997.25
$ws.Range("L91").Value = 1997.25
$ws.Range("N91").Value = -4805.25

# Row 110
$ws.Range("H110").Value = 695
$ws.Range("I110").Value = 695
$ws.Range("K110").Value = 695
$ws.Range("M110").Value = 1350

# Row 122
$ws.Range("H122").Value = 7294.95
$ws.Range("I122").Value = 6874.9165
$ws.Range("K122").Value = 20624.7495
$ws.Range("M122").Value = -18174.7495

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3793.3
$ws.Range("I20").Value = 5101
$ws.Range("J20").Value = 1831.75
$ws.Range("K20").Value = 5101
$ws.Range("L20").Value = 1831.75
$ws.Range("M20").Value = -4854
$ws.Range("N20").Value = -2325.75

# Row 86
$ws.Range("H86").Value = 3995
$ws.Range("I86").Value = 3995
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3995
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -2872

# Row 89
$ws.Range("H89").Value = 3995
$ws.Range("I89").Value = 3995
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 19975
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -14359

# Row 107
$ws.Range("H107").Value = 714.375
$ws.Range("I107").Value = 906
$ws.Range("J107").Value = 395
$ws.Range("K107").Value = 906
$ws.Range("L107").Value = 395
$ws.Range("M107").Value = 1014
$ws.Range("N107").Value = -4235

# Row 134
$ws.Range("H134").Value = 1635.3529
$ws.Range("I134").Value = 1675.125
$ws.Range("J134").Value = 999
$ws.Range("K134").Value = 5025.375
$ws.Range("L134").Value = 2997
$ws.Range("M134").Value = -2490.375
$ws.Range("N134").Value = -8067

$ws = $wb.Worksheets.Item("CRP")
# Row 36
$ws.Range("H36").Value = 212.75
$ws.Range("I36").Value = 212.75
$ws.Range("K36").Value = 212.75
$ws.Range("M36").Value = 175.25

# Row 40
$ws.Range("H40").Value = 212.75
$ws.Range("I40").Value = 212.75
$ws.Range("K40").Value = 212.75
$ws.Range("M40").Value = -52.75

# Row 58
$ws.Range("H58").Value = 1878.2858
$ws.Range("I58").Value = 1109.1333
$ws.Range("K58").Value = 1109.1333
$ws.Range("M58").Value = -906.1333

# Row 60
$ws.Range("H60").Value = 19947.75
$ws.Range("I60").Value = 13270.154
$ws.Range("J60").Value = 48884
$ws.Range("K60").Value = 13270.154
$ws.Range("L60").Value = 48884
$ws.Range("M60").Value = -12759.154
$ws.Range("N60").Value = -49906

# Row 62
$ws.Range("H62").Value = 59661.715
$ws.Range("I62").Value = 2608.25
$ws.Range("J62").Value = 135733
$ws.Range("K62").Value = 2608.25
$ws.Range("L62").Value = 135733
$ws.Range("M62").Value = -1984.25
$ws.Range("N62").Value = -136981

# Row 65
$ws.Range("H65").Value = 59661.715
$ws.Range("I65").Value = 2608.25
$ws.Range("J65").Value = 135733
$ws.Range("K65").Value = 13041.25
$ws.Range("L65").Value = 678665
$ws.Range("M65").Value = -9921.25
$ws.Range("N65").Value = -684905

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = 0

# Row 136
$ws.Range("H136").Value = 1878.2858
$ws.Range("I136").Value = 1109.1333
$ws.Range("K136").Value = 3327.3999
$ws.Range("M136").Value = -777.3998999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 200690.4
$ws.Range("J2").Value = 1084
$ws.Range("L2").Value = 6504
$ws.Range("N2").Value = -6730

# Row 4
$ws.Range("H4").Value = 3471800.2
$ws.Range("I4").Value = 4673303.5
$ws.Range("J4").Value = 791.55554
$ws.Range("K4").Value = 14019910.5
$ws.Range("L4").Value = 2374.66662
$ws.Range("M4").Value = -14019798.5
$ws.Range("N4").Value = -2598.66662

# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = 0

# Row 69
$ws.Range("H69").Value = 2816.6667
$ws.Range("I69").Value = 300
$ws.Range("K69").Value = 900
$ws.Range("M69").Value = -89

# Row 72
$ws.Range("H72").Value = 2816.6667
$ws.Range("I72").Value = 300
$ws.Range("K72").Value = 2700
$ws.Range("M72").Value = 1356

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 32250
$ws.Range("J32").Value = 32250
$ws.Range("L32").Value = 32250
$ws.Range("N32").Value = -32842

# Row 52
$ws.Range("H52").Value = 39166.332
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 70
$ws.Range("H70").Value = 6115.615
$ws.Range("I70").Value = 5250.8335
$ws.Range("J70").Value = 6856.857
$ws.Range("K70").Value = 5250.8335
$ws.Range("L70").Value = 6856.857
$ws.Range("M70").Value = -4980.8335
$ws.Range("N70").Value = -7396.857

# Row 73
$ws.Range("H73").Value = 6115.615
$ws.Range("I73").Value = 5250.8335
$ws.Range("J73").Value = 6856.857
$ws.Range("K73").Value = 5250.8335
$ws.Range("L73").Value = 6856.857
$ws.Range("M73").Value = -4314.8335
$ws.Range("N73").Value = -8728.857

# Row 80
$ws.Range("H80").Value = 23158.166
$ws.Range("I80").Value = 9333
$ws.Range("J80").Value = 36983.332
$ws.Range("K80").Value = 9333
$ws.Range("L80").Value = 36983.332
$ws.Range("M80").Value = -8335
$ws.Range("N80").Value = -38979.332

# Row 83
$ws.Range("H83").Value = 23158.166
$ws.Range("I83").Value = 9333
$ws.Range("J83").Value = 36983.332
$ws.Range("K83").Value = 46665
$ws.Range("L83").Value = 184916.66
$ws.Range("M83").Value = -41673
$ws.Range("N83").Value = -194900.66

# Row 113
$ws.Range("H113").Value = 2932.6
$ws.Range("I113").Value = 2221
$ws.Range("K113").Value = 2221
$ws.Range("M113").Value = -51

# Row 122
$ws.Range("H122").Value = 63173.59
$ws.Range("I122").Value = 3543.9092
$ws.Range("J122").Value = 172494.67
$ws.Range("K122").Value = 10631.7276
$ws.Range("L122").Value = 517484.01
$ws.Range("M122").Value = -8181.7276
$ws.Range("N122").Value = -522384.01

$ws = $wb.Worksheets.Item("LTW")
# Row 34
$ws.Range("H34").Value = 3173.6667
$ws.Range("I34").Value = 3173.6667
$ws.Range("K34").Value = 3173.6667
$ws.Range("M34").Value = -3001.6667

# Row 46
$ws.Range("H46").Value = 4358
$ws.Range("I46").Value = 3180
$ws.Range("J46").Value = 4811.077
$ws.Range("K46").Value = 3180
$ws.Range("L46").Value = 4811.077
$ws.Range("M46").Value = -2992
$ws.Range("N46").Value = -5187.077

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7249.875
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 7499.8335
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 7499.8335
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -8747.833500000001

# Row 65
$ws.Range("H65").Value = 7249.875
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 7499.8335
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 37499.1675
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -43739.1675

# Row 104
$ws.Range("H104").Value = 15498.5
$ws.Range("J104").Value = 15498.5
$ws.Range("L104").Value = 15498.5
$ws.Range("N104").Value = -22486.5
